$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (groups_example_A) updates ---
$ws1.Range("C5").Value = 29
$ws1.Range("D5").Value = 8.17
$ws1.Range("E5").Value = 67.58
$ws1.Range("F5").Value = 15.96
$ws1.Range("H5").Value = 35.92
$ws1.Range("I5").Value = 58.65

$ws1.Range("C6").Value = 29
$ws1.Range("D6").Value = 8.9
$ws1.Range("E6").Value = 63.34
$ws1.Range("F6").Value = 20.34
$ws1.Range("H6").Value = 4.15
$ws1.Range("I6").Value = 67.9

$ws1.Range("C7").Value = 29
$ws1.Range("D7").Value = 6.25
$ws1.Range("E7").Value = 38.05
$ws1.Range("F7").Value = 17.74
$ws1.Range("H7").Value = 46.47
$ws1.Range("I7").Value = 42.31

$ws1.Range("C8").Value = 22
$ws1.Range("D8").Value = 10.98
$ws1.Range("E8").Value = 55.12
$ws1.Range("F8").Value = 18.06
$ws1.Range("H8").Value = 156.13
$ws1.Range("I8").Value = 76.46

$ws1.Range("C26").Value = 20
$ws1.Range("D26").Value = 7.04
$ws1.Range("E26").Value = 12.67
$ws1.Range("F26").Value = 14.34
$ws1.Range("H26").Value = 18.43
$ws1.Range("I26").Value = 70.07

$ws1.Range("C27").Value = 20
$ws1.Range("D27").Value = 8.94
$ws1.Range("E27").Value = 24.67
$ws1.Range("F27").Value = 19.52
$ws1.Range("H27").Value = 6.1
$ws1.Range("I27").Value = 44.18

$ws1.Range("C28").Value = 20
$ws1.Range("D28").Value = 5.33
$ws1.Range("E28").Value = 12.67
$ws1.Range("F28").Value = 16.17
$ws1.Range("H28").Value = 34.11
$ws1.Range("I28").Value = 31.24

$ws1.Range("C29").Value = 42
$ws1.Range("D29").Value = 1.1
$ws1.Range("E29").Value = 46.86
$ws1.Range("F29").Value = 9.42
$ws1.Range("H29").Value = 52.98
$ws1.Range("I29").Value = 26.62

$ws1.Range("C30").Value = 42
$ws1.Range("D30").Value = 2.06
$ws1.Range("F30").Value = 6.92
$ws1.Range("H30").Value = 34.55
$ws1.Range("I30").Value = 91.28

$ws1.Range("C31").Value = 42
$ws1.Range("D31").Value = 0.46
$ws1.Range("F31").Value = 5.4
$ws1.Range("H31").Value = 51.94
$ws1.Range("I31").Value = 27.28

# --- Sheet2 (groups_example_B) updates ---
$ws2.Range("B4").Value = 5
$ws2.Range("C4").Value = 7.18
$ws2.Range("D4").Value = 64
$ws2.Range("E4").Value = 24.69
$ws2.Range("G4").Value = 10.4
$ws2.Range("H4").Value = 75.44
$ws2.Range("I4").Value = 5
$ws2.Range("J4").Value = 9.12
$ws2.Range("K4").Value = 85.33
$ws2.Range("L4").Value = 29.68
$ws2.Range("N4").Value = 9.41
$ws2.Range("O4").Value = 60.68
$ws2.Range("P4").Value = 5
$ws2.Range("Q4").Value = 1.95
$ws2.Range("R4").Value = 64
$ws2.Range("S4").Value = 20.8
$ws2.Range("T4").Value = 178.47
$ws2.Range("U4").Value = 20.51

$ws2.Range("B5").Value = 29
$ws2.Range("C5").Value = 8.9
$ws2.Range("D5").Value = 63.34
$ws2.Range("E5").Value = 20.34
$ws2.Range("G5").Value = 4.15
$ws2.Range("H5").Value = 67.9
$ws2.Range("I5").Value = 29
$ws2.Range("J5").Value = 8.17
$ws2.Range("K5").Value = 67.58
$ws2.Range("L5").Value = 15.96
$ws2.Range("N5").Value = 35.92
$ws2.Range("O5").Value = 58.65
$ws2.Range("P5").Value = 29
$ws2.Range("Q5").Value = 6.25
$ws2.Range("R5").Value = 38.05
$ws2.Range("S5").Value = 17.74
$ws2.Range("T5").Value = 46.47
$ws2.Range("U5").Value = 42.31

$ws2.Range("B6").Value = 22
$ws2.Range("C6").Value = 19.54
$ws2.Range("D6").Value = 79.76
$ws2.Range("E6").Value = 30.26
$ws2.Range("G6").Value = 5.26
$ws2.Range("H6").Value = 84.03
$ws2.Range("I6").Value = 22
$ws2.Range("J6").Value = 10.98
$ws2.Range("K6").Value = 55.12
$ws2.Range("L6").Value = 18.06
$ws2.Range("N6").Value = 156.13
$ws2.Range("O6").Value = 76.46
$ws2.Range("P6").Value = 22
$ws2.Range("Q6").Value = 11
$ws2.Range("R6").Value = 86.5
$ws2.Range("S6").Value = 23.45
$ws2.Range("T6").Value = 26.56
$ws2.Range("U6").Value = 45.16

$ws2.Range("B7").Value = 14
$ws2.Range("C7").Value = 3.6
$ws2.Range("D7").Value = 72.14
$ws2.Range("E7").Value = 22.88
$ws2.Range("G7").Value = 4.59
$ws2.Range("H7").Value = 34.18
$ws2.Range("I7").Value = 14
$ws2.Range("J7").Value = 3.01
$ws2.Range("K7").Value = 61.71
$ws2.Range("L7").Value = 12.78
$ws2.Range("N7").Value = 12.28
$ws2.Range("O7").Value = 84.46
$ws2.Range("P7").Value = 14
$ws2.Range("Q7").Value = 2.69
$ws2.Range("R7").Value = 61.22
$ws2.Range("S7").Value = 17.36
$ws2.Range("T7").Value = 35.14
$ws2.Range("U7").Value = 26.63

$ws2.Range("B8").Value = 20
$ws2.Range("C8").Value = 9.52
$ws2.Range("D8").Value = 45.33
$ws2.Range("E8").Value = 31.84
$ws2.Range("G8").Value = 177.11
$ws2.Range("H8").Value = 37.21
$ws2.Range("I8").Value = 20
$ws2.Range("J8").Value = 15.07
$ws2.Range("K8").Value = 84.67
$ws2.Range("L8").Value = 27.03
$ws2.Range("N8").Value = 179.69
$ws2.Range("O8").Value = 77.75
$ws2.Range("P8").Value = 20
$ws2.Range("Q8").Value = 15.21
$ws2.Range("R8").Value = 78
$ws2.Range("S8").Value = 32.15
$ws2.Range("T8").Value = 17.96
$ws2.Range("U8").Value = 36.81

$ws2.Range("B9").Value = 15
$ws2.Range("C9").Value = 3.04
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 11.04
$ws2.Range("G9").Value = 1.69
$ws2.Range("H9").Value = 73.6
$ws2.Range("I9").Value = 15
$ws2.Range("J9").Value = 2.09
$ws2.Range("K9").Value = 0
$ws2.Range("L9").Value = 8.92
$ws2.Range("N9").Value = 14.21
$ws2.Range("O9").Value = 84.55
$ws2.Range("P9").Value = 15
$ws2.Range("Q9").Value = 1.9
$ws2.Range("R9").Value = 0
$ws2.Range("S9").Value = 10.48
$ws2.Range("T9").Value = 38.57
$ws2.Range("U9").Value = 37.83

$ws2.Range("B10").Value = 17
$ws2.Range("C10").Value = 4.08
$ws2.Range("D10").Value = 0
$ws2.Range("E10").Value = 13.62
$ws2.Range("G10").Value = 0.11
$ws2.Range("H10").Value = 62.52
$ws2.Range("I10").Value = 17
$ws2.Range("J10").Value = 4.09
$ws2.Range("K10").Value = 0
$ws2.Range("L10").Value = 11.27
$ws2.Range("N10").Value = 175.98
$ws2.Range("O10").Value = 93.3
$ws2.Range("P10").Value = 17
$ws2.Range("Q10").Value = 3.61
$ws2.Range("R10").Value = 0
$ws2.Range("S10").Value = 14.44
$ws2.Range("T10").Value = 25.48
$ws2.Range("U10").Value = 32.96

$ws2.Range("B11").Value = 17
$ws2.Range("C11").Value = 1.6
$ws2.Range("D11").Value = 0
$ws2.Range("E11").Value = 9.26
$ws2.Range("G11").Value = 176.28
$ws2.Range("H11").Value = 64.91
$ws2.Range("I11").Value = 17
$ws2.Range("J11").Value = 10.8
$ws2.Range("K11").Value = 66.51
$ws2.Range("L11").Value = 20.45
$ws2.Range("N11").Value = 151.43
$ws2.Range("O11").Value = 87.94
$ws2.Range("P11").Value = 17
$ws2.Range("Q11").Value = 4.77
$ws2.Range("R11").Value = 66.44
$ws2.Range("S11").Value = 20.21
$ws2.Range("T11").Value = 119.33
$ws2.Range("U11").Value = 32.72

$ws2.Range("B12").Value = 20
$ws2.Range("C12").Value = 8.94
$ws2.Range("D12").Value = 24.67
$ws2.Range("E12").Value = 19.52
$ws2.Range("G12").Value = 6.1
$ws2.Range("H12").Value = 44.18
$ws2.Range("I12").Value = 20
$ws2.Range("J12").Value = 7.04
$ws2.Range("K12").Value = 12.67
$ws2.Range("L12").Value = 14.34
$ws2.Range("N12").Value = 18.43
$ws2.Range("O12").Value = 70.07
$ws2.Range("P12").Value = 20
$ws2.Range("Q12").Value = 5.33
$ws2.Range("R12").Value = 12.67
$ws2.Range("S12").Value = 16.17
$ws2.Range("T12").Value = 34.11
$ws2.Range("U12").Value = 31.24

$ws2.Range("B13").Value = 42
$ws2.Range("C13").Value = 2.06
$ws2.Range("D13").Value = 0
$ws2.Range("E13").Value = 6.92
$ws2.Range("G13").Value = 34.55
$ws2.Range("H13").Value = 91.28
$ws2.Range("I13").Value = 42
$ws2.Range("J13").Value = 1.1
$ws2.Range("K13").Value = 46.86
$ws2.Range("L13").Value = 9.42
$ws2.Range("N13").Value = 52.98
$ws2.Range("O13").Value = 26.62
$ws2.Range("P13").Value = 42
$ws2.Range("Q13").Value = 0.46
$ws2.Range("R13").Value = 0
$ws2.Range("S13").Value = 5.4
$ws2.Range("T13").Value = 51.94
$ws2.Range("U13").Value = 27.28
